# Material Database.xlsx — "Completed code and plotting for VT 1"
#
# 1. Remove the now-unused "Inputs" sheet.
# 2. Clean up the Materials list:
#      - drop "Oak Timber"
#      - row 9 becomes the (previously empty) "Wood" material, now with full data
#      - rename "Rockwool Loft Roll" -> "Tiles" (its row/Conductivity value is unchanged)
#      - add a new row 14 with just a Conductivity (column D) value for the plotting step
# 3. Update the view state (active cell / tab) and calculation mode.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Materials")

# --- Remove the "Inputs" worksheet (no longer needed) ---
$inputs = $wb.Worksheets.Item("Inputs")
[void]$inputs.Delete()

# --- Row 9: was "Oak Timber" (blank data row) -> now "Wood" with full properties ---
$ws.Range("A9").Value = "Wood"
$ws.Range("B9").Value = 360
$ws.Range("C9").Value = 1720
$ws.Range("D9").Value = 0.11
$ws.Range("E9").Value = 0.9
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 0.35
$ws.Range("H9").Value = 0.65

# --- Row 13: "Rockwool Loft Roll" -> "Tiles" (data untouched) ---
$ws.Range("A13").Value = "Tiles"

# --- New row 14: Conductivity-only entry used by the VT1 plot ---
$ws.Range("D14").Value = 0.8
$ws.Range("D13").Copy()
[void]$ws.Range("D14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- View state: select the sheet / active cell used while finishing the work ---
[void]$ws.Activate()
[void]$ws.Range("E9").Select()

# --- Calculation mode: Automatic except for data tables ---
$excel.Calculation = 2
